# "added AMC & GME" -- in this workbook that means adding the Taipei-listed
# Taiwan Semi ticker, a Samsung ticker, and a new Broadcom/Avago row to the
# "Semiconductors" sheet, then leaving the "Hardware" sheet as the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Semiconductors")
$ws2 = $wb.Worksheets.Item("Hardware")

# Row 8: Taiwan Semi's local (Taipei) ticker, pushing the old Samsung/Micron
# rows down by one.
$ws1.Range("B8").Value = "Taiwan Semi"
$ws1.Range("C8").Value = "2330 TT"

# Row 9: Samsung, now with its Korean exchange ticker.
$ws1.Range("B9").Value = "Samsung"
$ws1.Range("C9").Value = "005930 KS"

# Row 10: Micron, unchanged, just shifted down.
$ws1.Range("B10").Value = "Micron"
$ws1.Range("C10").Value = "MU"

# Row 11: new Broadcom / Avago entry.
$ws1.Range("B11").Value = "Broadcom"
$ws1.Range("C11").Value = "Avago"

# Leave the selection on the last-edited cell of the Semiconductors sheet,
# then switch over to (and leave active) the Hardware sheet.
$ws1.Range("D11").Select() | Out-Null
$ws2.Activate() | Out-Null
